$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "245.21" },
    @{ Cell = "E2"; Value = "-0.55%" },
    @{ Cell = "D3"; Value = "28.60" },
    @{ Cell = "E3"; Value = "-3.87%" },
    @{ Cell = "D4"; Value = "5.250" },
    @{ Cell = "E4"; Value = "1.28%" },
    @{ Cell = "D5"; Value = "0.05699" },
    @{ Cell = "E5"; Value = "-0.37%" },
    @{ Cell = "D6"; Value = "6.619" },
    @{ Cell = "E6"; Value = "0.29%" },
    @{ Cell = "E7"; Value = "3.22%" },
    @{ Cell = "D8"; Value = "0.8501" },
    @{ Cell = "E8"; Value = "-0.68%" },
    @{ Cell = "D9"; Value = "0.8539" },
    @{ Cell = "E9"; Value = "-1.89%" },
    @{ Cell = "D10"; Value = "0.1366" },
    @{ Cell = "E10"; Value = "0.19%" },
    @{ Cell = "D11"; Value = "0.07079" },
    @{ Cell = "E11"; Value = "-0.10%" },
    @{ Cell = "D12"; Value = "0.03277" },
    @{ Cell = "E12"; Value = "-2.86%" },
    @{ Cell = "D13"; Value = "0.03150" },
    @{ Cell = "E13"; Value = "7.72%" },
    @{ Cell = "D14"; Value = "0.09202" },
    @{ Cell = "E14"; Value = "-1.97%" },
    @{ Cell = "D15"; Value = "0.001536" },
    @{ Cell = "E15"; Value = "1.71%" },
    @{ Cell = "D16"; Value = "0.0005956" },
    @{ Cell = "E16"; Value = "-94.21%" },
    @{ Cell = "D17"; Value = "0.005916" },
    @{ Cell = "E17"; Value = "-2.62%" },
    @{ Cell = "E18"; Value = "0.09%" },
    @{ Cell = "D19"; Value = "2.175" },
    @{ Cell = "E19"; Value = "-4.48%" },
    @{ Cell = "E20"; Value = "-0.61%" },
    @{ Cell = "E21"; Value = "-1.15%" },
    @{ Cell = "D22"; Value = "3.495" },
    @{ Cell = "E22"; Value = "0.79%" },
    @{ Cell = "D23"; Value = "0.04071" },
    @{ Cell = "E23"; Value = "-2.34%" },
    @{ Cell = "E24"; Value = "-0.06%" },
    @{ Cell = "D25"; Value = "0.001218" },
    @{ Cell = "E25"; Value = "-0.32%" },
    @{ Cell = "D26"; Value = "0.004142" },
    @{ Cell = "E26"; Value = "-17.50%" },
    @{ Cell = "E27"; Value = "-0.85%" },
    @{ Cell = "D28"; Value = "0.0001448" },
    @{ Cell = "D40"; Value = "0.03755" },
    @{ Cell = "E40"; Value = "0.25%" },
    @{ Cell = "D41"; Value = "0.1063" },
    @{ Cell = "E41"; Value = "-0.84%" },
    @{ Cell = "D42"; Value = "0.003718" },
    @{ Cell = "E42"; Value = "-35.45%" },
    @{ Cell = "D43"; Value = "0.002299" },
    @{ Cell = "E43"; Value = "14.96%" },
    @{ Cell = "D44"; Value = "0.009327" },
    @{ Cell = "E44"; Value = "12.01%" },
    @{ Cell = "D45"; Value = "0.00005265" },
    @{ Cell = "E45"; Value = "1.00%" },
    @{ Cell = "E46"; Value = "-0.04%" },
    @{ Cell = "D47"; Value = "0.07495" },
    @{ Cell = "E47"; Value = "15.85%" },
    @{ Cell = "D48"; Value = "0.002438" },
    @{ Cell = "E48"; Value = "-2.97%" },
    @{ Cell = "D49"; Value = "0.00002099" },
    @{ Cell = "E49"; Value = "-0.04%" },
    @{ Cell = "D50"; Value = "0.0001999" },
    @{ Cell = "E50"; Value = "-0.04%" }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.Value = "'" + $u.Value
    $cell.Style = "Normal"
}
